$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" column (E16:E21) from descending (2506..2501)
# to ascending (2501..2506), and move the matching "Valor Mora" (F column)
# figures for periods 2501 / 2506 so they travel with their period.
$ws.Range("E16").Value = "2501"
$ws.Range("E17").Value = "2502"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2505"
$ws.Range("E21").Value = "2506"

$ws.Range("F16").Value = 36000
$ws.Range("F21").Value = 60000
